# Sincronização de dados: duas novas avaliações foram adicionadas,
# uma antes e outra depois da avaliação "Colaborador muito atencioso."
# (linha 19 -> agora linha 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere uma nova linha antes da antiga linha 19 (empurra tudo uma linha para baixo)
$ws.Rows.Item(19).Insert()

# Nova linha 19: avaliação sem comentário
$ws.Range("A19").Value = 5
$ws.Range("C19").Value = 45940.6359195949
$ws.Range("C19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D19").Value = "ZDljNWUxNGUtNDI1OS00ZTZhLWEzZmMtODlmMmZkMzNlMzBjOjU3MDE2"

# Nova linha 21 (acrescentada ao final): avaliação sem comentário
$ws.Range("A21").Value = 5
$ws.Range("C21").Value = 45940.63159275463
$ws.Range("C21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D21").Value = "NzYxOTRkMzgtZDQwNy00ODM5LWI5NDctYzJkMTczZDkzZGQzOjU3MDE2"
